# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to the refreshed values from the data pull.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"   = @{ 4 = 79; 5 = 6; 7 = 7632; 9 = 199; 10 = 1078; 11 = 640; 12 = 10; 17 = 747 }
    "全部类型" = @{ 4 = 79; 5 = 6; 8 = 7632; 10 = 199; 11 = 1078; 12 = 640; 13 = 10; 18 = 747 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
